$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.911.40'
$ws.Range('D3').Value = '1.636.54'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '215.53'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.67%  '
$ws.Range('E6').Value = '  +0.15%  '
$ws.Range('E7').Value = '  +0.08%  '
$ws.Range('E8').Value = '  -0.94%  '
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.53'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('B12').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C12').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D12').Value = '1.863.77'
$ws.Range('E12').Value = '  -0.62%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.27'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').Value = '1.637.75'
$ws.Range('E14').Value = '  -1.15%  '
$ws.Range('E15').Value = '  -0.59%  '
$ws.Range('D16').Value = '0.0₃0763'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '62.84'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.85%  '
$ws.Range('D18').Value = '25.932.16'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('E19').Value = '  +0.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '192.87'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.38%  '
$ws.Range('E21').Value = '  -1.98%  '
$ws.Range('E22').Value = '  -1.71%  '
$ws.Range('E23').Value = '  -0.79%  '
$ws.Range('E24').Value = '  +4.86%  '
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '143.11'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').Value = '  +0.09%  '
$ws.Range('E28').Value = '  -0.97%  '
$ws.Range('E29').Value = '  -0.62%  '
$ws.Range('E30').Value = '  -0.63%  '
$ws.Range('E31').Value = '  -0.67%  '
$ws.Range('E32').Value = '  -2.14%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.24'
$ws.Range('D33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.53'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.31%  '
$ws.Range('E35').Value = '  +1.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.900'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D37').Value = '1.134.03'
$ws.Range('E37').Value = '  -0.41%  '
$ws.Range('E38').Value = '  -1.62%  '
$ws.Range('E39').Value = '  -1.61%  '
$ws.Range('E40').Value = '  -0.69%  '
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.29'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.22%  '
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('D44').Value = '1.773.34'
$ws.Range('E44').Value = '  -0.63%  '
$ws.Range('D45').Value = '0.0₆0115'
$ws.Range('E45').Value = '  +2.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.50'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.05%  '
$ws.Range('E47').Value = '  +2.15%  '
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.69'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.30%  '
$ws.Range('E50').Value = '  -0.96%  '
$ws.Range('E51').Value = '  -1.28%  '
